# Add a new "report_order" stored-procedure row into the db-objects list.
# This mirrors the author's edit: a new row is inserted at row 129 (pushing
# the existing "view" rows down by one), and the sheet selection is moved
# to reflect where the user's cursor ended up after typing the new entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 129 ("view" / "view_ingredient" / ...)
# so the three existing "view" rows shift down to 130-132.
$ws.Rows.Item(129).Insert() | Out-Null

# Populate the newly inserted row with the new stored procedure entry.
$ws.Cells.Item(129, 1).Value2 = "stored procedure"
$ws.Cells.Item(129, 2).Value2 = "report_order"
$ws.Cells.Item(129, 3).Value2 = "input an order_id, product_id, lot_status, shows lot information"

# Move the selection to match the post-edit cursor position.
$ws.Range("C137").Select() | Out-Null
